# edit.ps1
# Applies the changes described by the commit diff to the active document.
#
# Summary of changes:
#  1. Delete the list paragraph "Ask users about the progress regarding their budget plans"
#  2. Remove the lastRenderedPageBreak marker before "View and add suggestions to suggestion board"
#  3. After the "...track each individual expense during their vacation." paragraph, insert a new
#     blank paragraph followed by a new paragraph of text about the web application's target audience
#  4. "User sign up page" -> "Sign up page"
#  5. "User login page" -> "Login page"
#  6. "User account page" -> "Profile page"
#  7. "User home page" -> "Home page"
#  8. Insert a new list item "View budget plan page" after "Add budget plan page"
#  9. Add a lastRenderedPageBreak marker to the "Develop frontend " run inside the schedule table

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Delete the "Ask users about the progress regarding their budget plans"
#    paragraph completely (including its paragraph mark).
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Ask users about the progress regarding their budget plans")
if ($found) {
    $rng.Expand(4) | Out-Null   # wdParagraph = 4 ; include the paragraph mark
    $rng.Delete()
}

# ---------------------------------------------------------------------------
# 2. Remove the lastRenderedPageBreak marker that precedes
#    "View and add suggestions to suggestion board".
#    The marker occupies a single, zero-width character position right
#    before the run's text, so re-assigning .Text on a range that starts
#    one character earlier drops the marker while preserving the text.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("View and add suggestions to suggestion board")
if ($found) {
    $full = $d.Range($rng.Start - 1, $rng.End)
    $full.Text = "View and add suggestions to suggestion board"
}

# ---------------------------------------------------------------------------
# 3. Insert two new paragraphs after the "...track each individual expense
#    during their vacation." paragraph: a blank paragraph, then a paragraph
#    of new text.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("track each individual expense during their vacation.")
if ($found) {
    $rng.Collapse(0) | Out-Null   # wdCollapseEnd = 0
    $newText = "This web application is for individual" + [char]0x2019 + "s above the age of 18 who have the means to fund their own vacations but require an estimated budget for their upcoming vacations. "
    $rng.InsertAfter("`r`r" + $newText)
}

# ---------------------------------------------------------------------------
# 4. "User sign up page" -> "Sign up page"
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("User sign up page")
if ($found) {
    $rng.Text = "Sign up page"
}

# ---------------------------------------------------------------------------
# 5. "User login page" -> "Login page"
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("User login page")
if ($found) {
    $rng.Text = "Login page"
}

# ---------------------------------------------------------------------------
# 6. "User account page" -> "Profile page"
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("User account page")
if ($found) {
    $rng.Text = "Profile page"
}

# ---------------------------------------------------------------------------
# 7. "User home page" -> "Home page"
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("User home page")
if ($found) {
    $rng.Text = "Home page"
}

# ---------------------------------------------------------------------------
# 8. Insert a new list item "View budget plan page" after "Add budget plan page"
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Add budget plan page")
if ($found) {
    $rng.Collapse(0) | Out-Null   # wdCollapseEnd = 0
    $rng.InsertAfter("`rView budget plan page")
}


